$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(1.185992460493708,0.2018160463589993,0,0.2228484590543829,1.769493182683249,0.002435500124603581,0,0.7229978485463384,0.03134599805285987,0,0.5403205836999803,0,1.215414870328779,2.768054110926158),
    @(1.089623443391474,0.1916701541959895,0,0.223208674916858,1.765991188807732,0.002438205388448357,0,0.7294918815187863,0.03054146774887556,0,0.5286973746315198,0,1.223174187351205,2.77300894922729),
    @(1.030617253160074,0.1853848362532631,0,0.2235024210391394,1.764805780481467,0.002439956057839818,0,0.7338986525173006,0.03004670748128646,0,0.5217619443411508,0,1.228402972053473,2.777827604818015),
    @(1.006614844690176,0.1828096328717947,0,0.2236404153005598,1.764565442499446,0.002440692074447609,0,0.7357998420550516,0.02984491084898977,0,0.5189865713680462,0,1.23065070315684,2.780237594672542),
    @(1.002631914665017,0.1823811872357055,0,0.2236644350539407,1.764540197880066,0.002440815656518748,0,0.7361218966840823,0.02981139243749453,0,0.5185288029784516,0,1.231031005671774,2.780664720301445),
    @(1.030293371603307,0.1853501622200184,0,0.2235042079658474,1.764801556266562,0.002439965892526898,0,0.7339238659782019,0.03004398667859576,0,0.5217243083678511,0,1.228432811980596,2.777858299981744),
    @(1.152731187773952,0.1983293725892423,0,0.2229576237235555,1.768085505455261,0.002436414337999345,0,0.7251499011652882,0.03106876700744365,0,0.5362712801052254,0,1.21799395645639,2.769393647389194),
    @(1.394078060659638,0.2233354278049546,0,0.2224599033821306,1.782178690763473,0.002430157835330128,0,0.7112752933673541,0.03307154319162819,0,0.5663852243118725,0,1.201203126200376,2.766907636500292),
    @(1.572091891930029,0.2414314119836547,0,0.2224421351436447,1.797199893658757,0.002425988559639066,0,0.7031171650943548,0.03453807366788197,0,0.5894675297765417,0,1.191102509413597,2.773715526586045),
    @(1.653213790067866,0.2496030770558946,0,0.2225091654238192,1.805047346032993,0.002424183736347953,0,0.6998488852668885,0.03520401774877158,0,0.6001742262661338,0,1.186991424173286,2.778694398062385),
    @(1.683951595097881,0.2526887040189649,0,0.222545310475649,1.808164780591838,0.002423513428689312,0,0.6986750477728876,0.03545600661263393,0,0.6042580336418553,0,1.185504109590383,2.780850856539558),
    @(1.677330855470132,0.2520245528175167,0,0.2225370479639039,1.807486903051071,0.002423657207928239,0,0.698925016023523,0.0354017450399482,0,0.6033772087997704,0,1.185821341443869,2.780374361672671),
    @(1.655742243197381,0.2498571106898169,0,0.2225119236329007,1.80530089811765,0.002424128326751946,0,0.699751033938675,0.03522475293410565,0,0.600509615618364,0,1.186867670277081,2.778866376029328),
    @(1.642520977403649,0.2485283390753636,0,0.2224979346458191,1.80398088902416,0.002424418610016894,0,0.7002653040058462,0.03511631489276112,0,0.5987569544659124,0,1.187517620405409,2.777978005617967),
    @(1.566793099694848,0.2408961496707036,0,0.2224392626866241,1.796707447114116,0.00242610835107173,0,0.7033396757043064,0.03449452713729073,0,0.5887719528241888,0,1.191380904149852,2.773428053244629),
    @(1.520371725443738,0.2361985105055453,0,0.2224224798152541,1.792505145277019,0.002427168421304415,0,0.7053392111452155,0.03411276279211251,0,0.5826991583610095,0,1.193874729804492,2.771119131734736),
    @(1.493684899582092,0.2334908861226666,0,0.2224198990855868,1.790183560995942,0.002427786790019949,0,0.7065309722175535,0.03389307134806785,0,0.5792256949687982,0,1.195354649334199,2.769968210381222),
    @(1.484651575724286,0.2325731619326916,0,0.2224202411656115,1.789413911408289,0.002427997645652661,0,0.7069416377268176,0.03381866916362597,0,0.5780529883404597,0,1.195863548698767,2.769608933170218),
    @(1.525311973608154,0.2366991696702598,0,0.2224235347038359,1.792942607347939,0.002427054680755805,0,0.7051220426018219,0.03415341385089121,0,0.5833436069415399,0,1.193604545612999,2.771346586967866),
    @(1.662082848599141,0.2504939810356746,0,0.2225190114872859,1.805939025042292,0.002423989591559096,0,0.6995066804256709,0.03527674510504397,0,0.601351101133119,0,1.186558453642448,2.779301947692204),
    @(1.751578707860745,0.2594583012410681,0,0.2226441230099816,1.815282588000997,0.002422062939039792,0,0.6962085756855032,0.03600979249187475,0,0.6132913594329921,0,1.182358274238794,2.786081434164259),
    @(1.703803775012659,0.2546786214631425,0,0.2225716233681183,1.810218032953387,0.002423084244465027,0,0.6979347775820983,0.035618659525678,0,0.6069030329337579,0,1.184562975948026,2.782318358494905),
    @(1.523078482786843,0.2364728431116134,0,0.22242303577017,1.79274453676851,0.002427106074993928,0,0.7052200930513877,0.03413503615985292,0,0.5830521960653527,0,1.193726552010254,2.771243204672828),
    @(1.328660474540527,0.2166188030955709,0,0.2225333025993272,1.777546633888889,0.002431775025048189,0,0.7146716270864033,0.03253054669856326,0,0.5580697114227888,0,1.205352360976008,2.76606609154075)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value2 = $rowVals[$c]
    }
}
